$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$mojibake = "$([char]194)$([char]177)"
$fixed = [char]177

for ($row = 2; $row -le 17; $row++) {
    foreach ($col in @("B", "C", "D")) {
        $cell = $ws.Range("$col$row")
        $val = $cell.Value()
        if ($val -ne $null) {
            $newVal = $val.Replace($mojibake, $fixed)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
